# Rename the original sheet and append a new, empty "Sheet1" after it,
# matching the commit "Added Sql test cases" (scaffolding a second sheet
# for SQL test cases while retitling the dashboard sheet).

$wb = $excel.ActiveWorkbook

# Rename the existing "TEST CASE TEMPLATE" sheet to "TC_Dashboard".
$dashboard = $wb.Worksheets.Item(1)
$dashboard.Name = "TC_Dashboard"

# Add a brand-new blank worksheet right after the dashboard sheet and name it "Sheet1".
$sheet1 = $wb.Worksheets.Add($null, $dashboard)
$sheet1.Name = "Sheet1"

# Keep the dashboard sheet active/selected, matching its prior tab state.
$dashboard.Activate()

# Restore the scrolled viewport (top-left visible cell moved from B15 to B14)
# while keeping the current selection (G30) intact.
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 2
